# Survey workbook update:
#  - Insert a new question "How many children..." right before the
#    zipcode question (becomes the new row 8).
#  - Insert a new question about NHL streaming on Hulu right before the
#    "sensitive nature" question (near the end of the sheet).
#  - Both new answer cells copy the borderless "plain" number format
#    already used by the streaming-services answer, matching how that
#    look is applied elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the "children in household" question above the zipcode row ---
$ws.Rows(8).Insert()
$ws.Range("A8").Value = "How many children under the age of 18 are living in your household? Please reference only the children for which you are the parent or legal guardian. (If there are no children under 18 in your household, please type 0)"
$ws.Range("B8").Value = 0
$ws.Rows(8).RowHeight = 18.75

# Borrow the plain (no border) number style already used by the
# streaming-services answer (now shifted down to row 12) for the new
# answer cell.
$ws.Range("B12").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = 0

# --- Insert the Hulu/NHL streaming frequency question before the
#     sensitive-nature closing question ---
$ws.Rows(13).Insert()
$ws.Range("A13").Value = "Please tell us how often you typically watch the NHL streaming on Hulu on your tablet, smartphone, smart TV or streaming device (Roku, Apple TV, Amazon TV Fire stick, etc.)."
$ws.Range("B13").Value = "rrrrrrrrrr"
$ws.Rows(13).RowHeight = 18.75

$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "rrrrrrrrrr"

# Streaming-services answer row height tightens slightly in the final
# layout.
$ws.Rows(12).RowHeight = 18.75

$excel.CutCopyMode = $false
